$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.438
$ws.Range("D4").Value = -7.88
$ws.Range("C7").Value = -12.995
$ws.Range("D12").Value = -7.4
$ws.Range("C16").Value = -12.856
$ws.Range("D18").Value = -8.489999999999998
$ws.Range("D19").Value = -8.004999999999999
$ws.Range("D20").Value = -7.683999999999999
$ws.Range("C28").Value = -12.798
$ws.Range("C29").Value = -12.266
$ws.Range("D31").Value = -7.973000000000001
$ws.Range("C32").Value = -12.475
$ws.Range("C40").Value = -12.173
$ws.Range("D40").Value = -7.587999999999999
$ws.Range("D42").Value = -8.111000000000001
$ws.Range("D47").Value = -7.475
$ws.Range("D48").Value = -7.404000000000001
$ws.Range("C52").Value = -11.259
$ws.Range("C57").Value = -13.742
$ws.Range("D63").Value = -6.987
$ws.Range("D64").Value = -7.640000000000001
$ws.Range("C66").Value = -11.491
$ws.Range("D76").Value = -7.649999999999999
$ws.Range("D81").Value = -7.851000000000001
$ws.Range("D89").Value = -8.217000000000001
$ws.Range("D94").Value = -7.683000000000002
$ws.Range("C100").Value = -11.54
